$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E4").Value = "12,5%"
$ws.Range("F4").Value = "87,5%"
$ws.Range("E5").Value = "12,5%"
$ws.Range("F5").Value = "87,5%"
$ws.Range("D7").Value = "12,5%"
$ws.Range("F7").Value = "87,5%"
$ws.Range("D8").Value = "12,5%"
$ws.Range("E8").Value = "12,5%"
$ws.Range("B9").Value = "12,5%"
$ws.Range("C9").Value = "12,5%"
$ws.Range("D9").Value = "12,5%"
$ws.Range("E9").Value = "12,5%"
$ws.Range("D11").Value = "12,5%"
$ws.Range("E11").Value = "12,5%"
$ws.Range("E12").Value = "12,5%"
$ws.Range("F12").Value = "87,5%"
$ws.Range("E13").Value = "12,5%"
$ws.Range("F13").Value = "62,5%"
